$wb = $excel.ActiveWorkbook

# Select the full data range on the "2024" sheet (as if copying last year's
# layout as a starting point / reference before building the new sheet).
$ws2024 = $wb.Worksheets.Item("2024")
$ws2024.Range("A1:C14").Select() | Out-Null

# Add the new "2025" sheet after the last existing sheet ("2024")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025"

# Copy the region-label formatting (A2:A14) from "2024" and paste it into the
# new sheet so the labels keep the same style used on every other year tab.
$ws2024.Range("A2:A14").Copy() | Out-Null
$ws.Range("A2:A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row
$ws.Range("A1").Value = "Región"
$ws.Range("B1").Value = "Familiar"
$ws.Range("C1").Value = "No Familiar"

# Region labels
$ws.Range("A2").Value = "San Juan"
$ws.Range("A3").Value = "Arecibo"
$ws.Range("A4").Value = "Ponce"
$ws.Range("A5").Value = "Humacao"
$ws.Range("A6").Value = "Mayaguez"
$ws.Range("A7").Value = "Caguas"
$ws.Range("A8").Value = "Bayamon"
$ws.Range("A9").Value = "Carolina"
$ws.Range("A10").Value = "Guayama"
$ws.Range("A11").Value = "Aguadilla"
$ws.Range("A12").Value = "Utuado"
$ws.Range("A13").Value = "Fajardo"
$ws.Range("A14").Value = "Aibonito"

# Familiar / No Familiar counts for 2025
$ws.Range("B2").Value = 132
$ws.Range("C2").Value = 154
$ws.Range("B3").Value = 75
$ws.Range("C3").Value = 57
$ws.Range("B4").Value = 87
$ws.Range("C4").Value = 68
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 34
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 41
$ws.Range("B7").Value = 105
$ws.Range("C7").Value = 91
$ws.Range("B8").Value = 122
$ws.Range("C8").Value = 191
$ws.Range("B9").Value = 74
$ws.Range("C9").Value = 85
$ws.Range("B10").Value = 43
$ws.Range("C10").Value = 31
$ws.Range("B11").Value = 65
$ws.Range("C11").Value = 39
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 27
$ws.Range("B13").Value = 21
$ws.Range("C13").Value = 21
$ws.Range("B14").Value = 31
$ws.Range("C14").Value = 32

# Widen the "No Familiar" column so the header fits, matching the other
# year tabs' column-C treatment.
$ws.Columns.Item(3).ColumnWidth = 11

# Leave the same selection state a user pasting the counts into B2:C14 would
# end up with, then make the new sheet the active tab.
$ws.Range("B2:C14").Select() | Out-Null
$ws.Activate() | Out-Null
